# chore: adapt column header formatting to respective input file names
#
# 1. Rename the header row from the generic "_old"/"_new" suffixes to the
#    concrete format-version suffixes "_FV2310"/"_FV2404".
# 2. Freeze the header row (pane split under row 1).
# 3. Turn the used range into a native Excel Table ("Table1") so the new
#    headers double as table column names.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310",
    "diff",
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- Freeze the header row --------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- Turn A1:U89 into a real table -----------------------------------------
$usedRange = $ws.Range("A1:U89")
$tbl = $ws.ListObjects.Add(1, $usedRange, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

$ws.Range("A1").Select()

Write-Host "done"
